# "contingencies with rene fine"
# Insert two new line entries (line7, line8) into the lines/extr table on
# Sheet1. This pushes the existing extr1-extr8 rows (currently rows 8-15)
# down to rows 10-17, and a handful of from_bus/to_bus/in_service values
# get refreshed with the source data's new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: push rows 8-15 (extr1..extr8, with full formatting) down to
#     rows 10-17.
$ws.Range("A8:E15").Copy($ws.Range("A10"))

# --- Step 2: the copy brought the old row index (A8:A15 = 6..13) along
#     with it, so re-number column A for the relocated rows (8..15), and
#     refresh from_bus/to_bus/in_service -- only extr1 (row10) and extr2
#     (row11) actually change (in_service flips false -> true); the rest
#     keep their copied values.
for ($i = 0; $i -le 7; $i++) {
    $ws.Cells.Item(10 + $i, 1).Value = 8 + $i
}
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(11, 5).Value = $true

# --- Step 3: overwrite rows 8-9 (which already carry the correct index
#     style from the source rows) with the two new line entries.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $false

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true
